$p = $ppt.ActivePresentation
for ($si = 1; $si -le $p.Slides.Count; $si++) {
    $s = $p.Slides.Item($si)
    Write-Output "SLIDE $si"
    for ($i = 1; $i -le $s.Shapes.Count; $i++) {
        $sh = $s.Shapes.Item($i)
        Write-Output "  $i : $($sh.Name)"
        if ($sh.HasTextFrame) {
            Write-Output "      text: $($sh.TextFrame.TextRange.Text)"
        }
    }
}
